$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, pushing existing rows 66-122 down to 67-123.
$ws.Rows.Item(66).Insert()

# Populate the new row 66 with the new weekly data point (dated 2022-08-17).
$ws.Cells.Item(66, 1).Value() = 11
$ws.Cells.Item(66, 2).Value() = "Vega Monumental Concepción"
$ws.Cells.Item(66, 3).Value() = "Bíobío"
$ws.Cells.Item(66, 4).Value() = 44790
$ws.Cells.Item(66, 5).Value() = 8
$ws.Cells.Item(66, 6).Value() = "Fruta"
$ws.Cells.Item(66, 7).Value() = 100108
$ws.Cells.Item(66, 8).Value() = "Tropicales y subtropicales"
$ws.Cells.Item(66, 9).Value() = 100108002
$ws.Cells.Item(66, 10).Value() = "Mango"
$ws.Cells.Item(66, 11).Value() = "Sin especificar"
$ws.Cells.Item(66, 12).Value() = "Primera"
$ws.Cells.Item(66, 13).Value() = 200
$ws.Cells.Item(66, 14).Value() = 9500
$ws.Cells.Item(66, 15).Value() = 10000
$ws.Cells.Item(66, 16).Value() = 9750
$ws.Cells.Item(66, 17).Value() = "$/bandeja 4 kilos"
$ws.Cells.Item(66, 18).Value() = "Brasil"
$ws.Cells.Item(66, 19).Value() = 2438
$ws.Cells.Item(66, 20).Value() = 4

# Make sure the date cell keeps the existing date number format used by column D.
$ws.Cells.Item(66, 4).NumberFormat = $ws.Cells.Item(67, 4).NumberFormat
